$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Milestone Payments")

$ws.Range("A1").Value = "SOW: SOW-1437 - fo check"
$ws.Range("A3").Value = "Total Contract Value: `$10.00"
$ws.Range("E9").Value = 3
$ws.Range("E11").Value = 3
